$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the previous numeric-format style/width from column B (was used for the old
# "cost" column) before repurposing it as a text ("Input") column.
$ws.Columns.Item(2).ClearFormats()

# Fill the new columns. The order below (column-wise) matches the order in which the
# new strings were first entered by the author, producing the same shared-string table.
$ws.Range("B1:B4").Value = "Input"
$ws.Range("E1:E4").Value = "Output"
$ws.Range("D1:D4").Value = "1 Mil Calls"
$ws.Range("G1:G4").Value = "1 Mil Calls"

for ($r = 1; $r -le 4; $r++) {
  $ws.Cells.Item($r, 3).Value = $r
  $ws.Cells.Item($r, 6).Value = $r + 4
}

# Give column A (the provider-name column) an explicit custom width.
$ws.Columns.Item(1).ColumnWidth = 20.5

# Update the active selection to reflect the new used range.
$ws.Range("G4").Select()
